# Apply updated loading_percent values for Case_1_115 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updatedValues = @{
    "B2" = 13.10499849755848
    "D2" = 10.51806254340979
    "E2" = 16.42452409154901
    "F2" = 32.34918952604533
    "G2" = 32.28741231828452
    "H2" = 15.31626555800419
    "I2" = 27.85111940153915
    "J2" = 11.70334645885426
    "K2" = 9.540977062126533
    "L2" = 9.215304311342999
    "M2" = 14.19114561851442
    "O2" = 23.75543092574356
    "B3" = 12.96293748519843
    "D3" = 10.53040863978117
    "E3" = 16.46823739214567
    "F3" = 32.45454654158664
    "G3" = 32.40633472478095
    "H3" = 15.36501161168746
    "I3" = 27.9731926771347
    "J3" = 11.72337599225397
    "K3" = 9.24213515442953
    "L3" = 9.184036116095843
    "M3" = 14.14543502211589
    "O3" = 23.84104362915687
    "B4" = 12.87649645396858
    "D4" = 10.53924792734677
    "E4" = 16.49663856551001
    "F4" = 32.52556995051417
    "G4" = 32.48785409739483
    "H4" = 15.39703925733887
    "I4" = 28.05225745434139
    "J4" = 11.7363467164755
    "K4" = 9.052245619194609
    "L4" = 9.165637381799254
    "M4" = 14.11864349019054
    "O4" = 23.89790533803196
    "B5" = 12.84150238091103
    "D5" = 10.54316711637003
    "E5" = 16.50860578016531
    "F5" = 32.55610384446689
    "G5" = 32.52320545411406
    "H5" = 15.4106186874034
    "I5" = 28.08551300421956
    "J5" = 11.74180197727847
    "K5" = 8.97333250836834
    "L5" = 9.158345349380351
    "M5" = 14.10805287042042
    "O5" = 23.92215656022656
    "B6" = 12.83570659966947
    "D6" = 10.54383706378372
    "E6" = 16.5106167224776
    "F6" = 32.56127004101775
    "G6" = 32.52920408039734
    "H6" = 15.41290544047409
    "I6" = 28.09109770028195
    "J6" = 11.74271807577554
    "K6" = 8.960138954675173
    "L6" = 9.157147056446368
    "M6" = 14.10631425110065
    "O6" = 23.92624865032089
    "B7" = 12.87602353035072
    "D7" = 10.5392994981822
    "E7" = 16.49679836495911
    "F7" = 32.52597530016724
    "G7" = 32.48832223568326
    "H7" = 15.39722025614308
    "I7" = 28.05270175261601
    "J7" = 11.73641960074842
    "K7" = 9.051187462789189
    "L7" = 9.165538200721008
    "M7" = 14.1184993283283
    "O7" = 23.89822802776103
    "B8" = 13.05587339750793
    "D8" = 10.52205858953592
    "E8" = 16.43927303161896
    "F8" = 32.38420123442097
    "G8" = 32.32664896545646
    "H8" = 15.33263819076467
    "I8" = 27.89235829030628
    "J8" = 11.71011335032677
    "K8" = 9.439308856200324
    "L8" = 9.20435968371824
    "M8" = 14.17512398921978
    "O8" = 23.78405846581781
    "B9" = 13.4132846952546
    "D9" = 10.49821008037186
    "E9" = 16.3388078664165
    "F9" = 32.1564982038633
    "G9" = 32.07731220765686
    "H9" = 15.22261194073221
    "I9" = 27.61044586453852
    "J9" = 11.66384115160151
    "K9" = 10.14659672280114
    "L9" = 9.286636409094994
    "M9" = 14.29598404475377
    "O9" = 23.59427162499929
    "B10" = 13.6767416976587
    "D10" = 10.48672400978846
    "E10" = 16.27245829683824
    "F10" = 32.01995585183504
    "G10" = 31.93573168774281
    "H10" = 15.15187558910333
    "I10" = 27.42300862798987
    "J10" = 11.63305419656798
    "K10" = 10.62985144544308
    "L10" = 9.350556286864265
    "M10" = 14.39035517978081
    "O10" = 23.47564845109961
    "B11" = 13.79635940405449
    "D11" = 10.48280062773545
    "E11" = 16.2438814293013
    "F11" = 31.96453119426475
    "G11" = 31.88042233339938
    "H11" = 15.12188237887412
    "I11" = 27.34198290987657
    "J11" = 11.61973879368934
    "K11" = 10.84116756671125
    "L11" = 9.380326851396777
    "M11" = 14.43440551350363
    "O11" = 23.42620845469214
    "B12" = 13.84158748749141
    "D12" = 10.48150132623885
    "E12" = 16.23329005466547
    "F12" = 31.94450628466766
    "G12" = 31.86079091579351
    "H12" = 15.11083849036051
    "I12" = 27.31190818266288
    "J12" = 11.61479529324858
    "K12" = 10.91991790602098
    "L12" = 9.391694370010097
    "M12" = 14.45123917492089
    "O12" = 23.40813761617123
    "B13" = 13.83185041594036
    "D13" = 10.48177287689054
    "E13" = 16.23556087808998
    "F13" = 31.94877614884349
    "G13" = 31.86496042006507
    "H13" = 15.11320303415883
    "I13" = 27.31835830287835
    "J13" = 11.61585557922431
    "K13" = 10.90301478402568
    "L13" = 9.389242079136695
    "M13" = 14.44760708595296
    "O13" = 23.41200052926323
    "B14" = 13.80008191998854
    "D14" = 10.48269000236234
    "E14" = 16.24300546458193
    "F14" = 31.96286442174885
    "G14" = 31.87878090252433
    "H14" = 15.12096750255319
    "I14" = 27.33949646993333
    "J14" = 11.61933011220932
    "K14" = 10.84767210845398
    "L14" = 9.381260211627454
    "M14" = 14.43578741604407
    "O14" = 23.42470870456945
    "B15" = 13.78061283215029
    "D15" = 10.48327601875411
    "E15" = 16.24759542465869
    "F15" = 31.9716193703294
    "G15" = 31.88741748425848
    "H15" = 15.12576433173036
    "I15" = 27.35252332790335
    "J15" = 11.62147121194407
    "K15" = 10.81360640852799
    "L15" = 9.376383167012191
    "M15" = 14.42856716899405
    "O15" = 23.43257763127436
    "B16" = 13.66891672183603
    "D16" = 10.48700653938565
    "E16" = 16.27435810317661
    "F16" = 32.02371270146471
    "G16" = 31.93952976153519
    "H16" = 15.15387965779196
    "I16" = 27.42838902150256
    "J16" = 11.63393823191038
    "K16" = 10.61586600576961
    "L16" = 9.348624180105128
    "M16" = 14.38749822243134
    "O16" = 23.4789705326021
    "B17" = 13.60031038971499
    "D17" = 10.48962794323302
    "E17" = 16.29118683084044
    "F17" = 32.05738449784007
    "G17" = 31.97383242349597
    "H17" = 15.17168692352655
    "I17" = 27.47601493925819
    "J17" = 11.64176269605384
    "K17" = 10.49234267803062
    "L17" = 9.331768502746783
    "M17" = 14.362584835263
    "O17" = 23.50858985420405
    "B18" = 13.56083025686831
    "D18" = 10.4912582761195
    "E18" = 16.30101747416367
    "F18" = 32.07738116747937
    "G18" = 31.99441832346491
    "H18" = 15.18213488375872
    "I18" = 27.50380734221508
    "J18" = 11.64632806766326
    "K18" = 10.42049584732448
    "L18" = 9.322139143964685
    "M18" = 14.34836128697111
    "O18" = 23.52605172107265
    "B19" = 13.54746068925949
    "D19" = 10.49183135361032
    "E19" = 16.30437195688482
    "F19" = 32.08425977873868
    "G19" = 32.00153520126142
    "H19" = 15.18570772147933
    "I19" = 27.51328599710285
    "J19" = 11.64788499188013
    "K19" = 10.39603390911926
    "L19" = 9.318890241198458
    "M19" = 14.3435638905262
    "O19" = 23.53203708711778
    "B20" = 13.60761593676907
    "D20" = 10.48933621025139
    "E20" = 16.28937974123646
    "F20" = 32.05373491238129
    "G20" = 31.97009222503227
    "H20" = 15.16977002373798
    "I20" = 27.47090377296198
    "J20" = 11.64092304986762
    "K20" = 10.50557503696384
    "L20" = 9.333556071271
    "M20" = 14.36522600000905
    "O20" = 23.50539277190285
    "B21" = 13.80941523664084
    "D21" = 10.48241556804486
    "E21" = 16.24081257197918
    "F21" = 31.95870020352593
    "G21" = 31.87468581864241
    "H21" = 15.11867837367099
    "I21" = 27.33327119246331
    "J21" = 11.61830688090675
    "K21" = 10.86396241035967
    "L21" = 9.383602172379307
    "M21" = 14.43925506139755
    "O21" = 23.42095833143518
    "B22" = 13.94088728859338
    "D22" = 10.47897860526986
    "E22" = 16.21041172121587
    "F22" = 31.90220434478407
    "G22" = 31.81998761434212
    "H22" = 15.08711646808041
    "I22" = 27.2468631478571
    "J22" = 11.60410132242885
    "K22" = 11.09076649597379
    "L22" = 9.416855714771057
    "M22" = 14.48852370751025
    "O22" = 23.36957047887456
    "B23" = 13.8707677164026
    "D23" = 10.48071387081184
    "E23" = 16.22651484382805
    "F23" = 31.93184309030726
    "G23" = 31.8484790616021
    "H23" = 15.10379436288962
    "I23" = 27.29265717087307
    "J23" = 11.61163058671309
    "K23" = 10.97040972915102
    "L23" = 9.399059642448691
    "M23" = 14.46214980920501
    "O23" = 23.3966496734862
    "B24" = 13.604313212855
    "D24" = 10.48946771876681
    "E24" = 16.29019624176396
    "F24" = 32.05538290180894
    "G24" = 31.97178047523373
    "H24" = 15.17063599895685
    "I24" = 27.47321324891972
    "J24" = 11.64130244523279
    "K24" = 10.49959527432975
    "L24" = 9.332747720502342
    "M24" = 14.36403161889647
    "O24" = 23.50683682309417
    "B25" = 13.31629793184839
    "D25" = 10.50359901627106
    "E25" = 16.36467149089358
    "F25" = 32.21270289331633
    "G25" = 32.13748055303582
    "H25" = 15.25060106608699
    "I25" = 27.68324385356134
    "J25" = 11.67579324507264
    "K25" = 9.961413652969313
    "L25" = 9.263748142039855
    "M25" = 14.26227842802447
    "O25" = 23.64196077511367
}

foreach ($cellRef in $updatedValues.Keys) {
    $ws.Range($cellRef).Value = $updatedValues[$cellRef]
}
